$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.241.48'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '2.006.81'
$ws.Range('E3').Value = '  -0.85%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '259.03'
$ws.Range('E5').Value = '  +4.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.609'
$ws.Range('E6').Value = '  -1.90%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.54'
$ws.Range('E8').Value = '  -4.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.378'
$ws.Range('E9').Value = '  -2.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0762'
$ws.Range('E10').Value = '  -5.11%  '
$ws.Range('E11').Value = '  -3.14%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').Value = '2.302.28'
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.18'
$ws.Range('E13').Value = '  -5.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.93'
$ws.Range('E14').Value = '  +2.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.772'
$ws.Range('E15').Value = '  -7.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.17'
$ws.Range('E16').Value = '  -3.92%  '
$ws.Range('D17').Value = '2.030.69'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').Value = '37.055.06'
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.77'
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').Value = '0.0₃0825'
$ws.Range('E20').Value = '  -3.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '234.60'
$ws.Range('E21').Value = '  +2.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.06'
$ws.Range('E22').Value = '  -3.15%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.56'
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.96'
$ws.Range('E26').Value = '  +0.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.83'
$ws.Range('E27').Value = '  -3.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.39'
$ws.Range('E28').Value = '  -2.31%  '
$ws.Range('E29').Value = '  -6.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.30'
$ws.Range('E30').Value = '  -4.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.118'
$ws.Range('E31').Value = '  -2.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.54'
$ws.Range('E32').Value = '  -4.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0620'
$ws.Range('E33').Value = '  -7.00%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.35'
$ws.Range('E34').Value = '  -4.65%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.40'
$ws.Range('E35').Value = '  -3.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.45'
$ws.Range('E36').Value = '  -2.41%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.80'
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('B38').Value = 'BinanceUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.33'
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.05'
$ws.Range('E40').Value = '  +1.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.17'
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('D42').Value = '1.442.65'
$ws.Range('E42').Value = '  +3.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0913'
$ws.Range('E43').Value = '  -5.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0207'
$ws.Range('E44').Value = '  -4.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '89.13'
$ws.Range('E45').Value = '  -1.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '15.52'
$ws.Range('E46').Value = '  -5.56%  '
$ws.Range('E47').Value = '  -1.51%  '
$ws.Range('E48').Value = '  +1.99%  '
$ws.Range('E49').Value = '  -8.89%  '
$ws.Range('D50').Value = '2.193.83'
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.92'
$ws.Range('E51').Value = '  -7.70%  '
